# The workbook has a typo in column A of rows 11-13: "ogistique DFT ClientIT"
# (missing leading "L"). Fix it to read "Logistique DFT ClientIT", matching the
# spelling used elsewhere in the sheet (e.g. "Logistique DSK ClientIT").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11:A13").Value = "Logistique DFT ClientIT"

# Leave the selection on the range that was just edited.
$ws.Range("A11:A13").Select()
